# Add a "Save" column (H) to the sheet, mirroring the style of the other
# header cells (B1:G1) and the plain numeric style of the data cells
# (B2:G14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same text + style as the rest of the header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new "Save" column, rows 2-14.
$saveValues = @(0, 0, 0, 1, 1, 0, 0, 0, 1, 0, 0, 1, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
